$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the old "Live Latency" column (D). This shifts chunk-size,
#    bitrate (kbps) and Segment Size (KB) one column to the left.
$ws.Columns("D").Delete()

# 2. The bitrate column (now E) used a flat 500 kbps for the first set of
#    rows; the new test run used 2000 kbps instead.
$ws.Range("E2:E7").Value = 2000

# 3. Add the new "Live Latency" results in column I for the first set of
#    rows (and the single extra sample in row 9).
$ws.Range("I2").Value = 6
$ws.Range("I3").Value = 6.5
$ws.Range("I4").Value = 5
$ws.Range("I5").Value = 6
$ws.Range("I6").Value = 6.5
$ws.Range("I7").Value = 6.5
$ws.Range("I9").Value = 6.5

# 4. Add the new "Not Chunked" column with header and values for the
#    first set of rows.
$ws.Range("L1").Value = "Not Chunked"
$ws.Range("L2").Value = 6.5
$ws.Range("L3").Value = 8
$ws.Range("L4").Value = 10
$ws.Range("L5").Value = 12
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 28

# 5. The header row no longer needs the taller wrapped height.
$ws.Rows("1").RowHeight = 30

# 6. Update the last selected cell.
$ws.Range("O4").Select() | Out-Null

# 7. Restore the window size/position recorded the last time the sheet
#    was saved.
$aw = $excel.ActiveWindow
$aw.Left = 3525
$aw.Top = 2310
$aw.Width = 21600
$aw.Height = 11835
